# Cập nhật data.xlsx từ công cụ QR
#
# The QR scanning tool appended a freshly scanned record. In the sheet
# that means: push the existing data row down by one and put the new
# record in its place, directly under the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the current data row (row 2) down to row 3 by inserting a new,
# blank row above it.
$ws.Rows.Item(2).Insert()

# Fill in the newly scanned record.
$ws.Range("A2").Value = "3ukwiw1n85x"
$ws.Range("B2").Value = "01mq60rp"
$ws.Range("C2").Value = "BAC"
$ws.Range("D2").Value = "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam"
$ws.Range("E2").Value = "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400"
$ws.Range("F2").Value = "2025-08-22T08:39:32.647Z"

# note / phone / branch / cccd / customerCode were not captured for this
# scan - write them as empty text (not left blank) so the row keeps a
# text value in every column, same as the other records.
foreach ($col in @("G", "H", "I", "J", "K")) {
    $ws.Range($col + "2").Value = "'"
}

$ws.Range("L2").Value = "Nguyễn Văn B"
$ws.Range("M2").Value = "68db92fc0a059591"
$ws.Range("N2").Value = "e0dd258af995cdce51b0cf5989fa8ea2457eed2fba41f7ab9fd82820e7514628"

# Clean up the quote-prefix formatting the empty-text trick above applied,
# so the new cells keep the same (default) formatting as the rest of the
# sheet.
$ws.Range("A1").Copy()
$ws.Range("G2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
